$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

$row = 19

$ws.Cells.Item($row, 1).Value = "13/06/2024 06:44:42"
$ws.Cells.Item($row, 2).Value = 1
$ws.Cells.Item($row, 3).Value = "APLAPOLLO"
$ws.Cells.Item($row, 4).Value = "Apl Apollo Tubes Limited"

# bsecode looks numeric but must be stored as text (matches source data).
# Use the leading-apostrophe trick to force text, then reset the style so
# no extra number-format style is attached to the cell itself.
$eCell = $ws.Cells.Item($row, 5)
$eCell.Value = "'533758"
$eCell.Style = "Normal"

$ws.Cells.Item($row, 6).Value = -1.7
$ws.Cells.Item($row, 7).Value = 1554.05
$ws.Cells.Item($row, 8).Value = 253346
